$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.410.97"
$ws.Range("E2").Value = "  +1.73%  "

$ws.Range("D3").Value = "2.178.78"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'253.50"
$ws.Range("E5").Value = "  +6.50%  "

$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").Value = "'75.38"
$ws.Range("E7").Value = "  +4.07%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "'0.586"
$ws.Range("E9").Value = "  +0.86%  "

$ws.Range("D10").Value = "'41.15"
$ws.Range("E10").Value = "  +3.02%  "

$ws.Range("D11").Value = "'0.0913"
$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'6.79"
$ws.Range("E12").Value = "  +1.17%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.101"
$ws.Range("E13").Value = "  +0.38%  "

$ws.Range("D14").Value = "2.505.69"
$ws.Range("E14").Value = "  +0.09%  "

$ws.Range("D15").Value = "'14.19"
$ws.Range("E15").Value = "  -1.16%  "

$ws.Range("D16").Value = "2.166.53"

$ws.Range("E17").Value = "  -1.13%  "

$ws.Range("D18").Value = "42.320.67"
$ws.Range("E18").Value = "  +1.79%  "

$ws.Range("E19").Value = "  -0.36%  "

$ws.Range("D20").Value = "'70.65"
$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("D21").Value = "'5.89"
$ws.Range("E21").Value = "  +1.58%  "

$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'9.63"
$ws.Range("E22").Value = "  -3.65%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'227.02"
$ws.Range("E23").Value = "  +0.44%  "

$ws.Range("E24").Value = "  +5.46%  "

$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("D26").Value = "'10.49"
$ws.Range("E26").Value = "  -1.85%  "

$ws.Range("D27").Value = "'3.34"
$ws.Range("E27").Value = "  +2.41%  "

$ws.Range("D28").Value = "'2.18"
$ws.Range("E28").Value = "  -0.57%  "

$ws.Range("D29").Value = "'37.27"
$ws.Range("E29").Value = "  +12.75%  "

$ws.Range("E30").Value = "  +2.88%  "

$ws.Range("D31").Value = "'169.38"
$ws.Range("E31").Value = "  -1.19%  "

$ws.Range("D32").Value = "'20.03"
$ws.Range("E32").Value = "  +0.81%  "

$ws.Range("D33").Value = "'0.0819"
$ws.Range("E33").Value = "  +5.89%  "

$ws.Range("D34").Value = "'5.13"
$ws.Range("E34").Value = "  -2.92%  "

$ws.Range("D35").Value = "'0.121"
$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("E36").Value = "  +4.75%  "

$ws.Range("D37").Value = "'4.25"
$ws.Range("E37").Value = "  -0.92%  "

$ws.Range("D38").Value = "'0.0335"
$ws.Range("E38").Value = "  +8.68%  "

$ws.Range("D39").Value = "'11.92"
$ws.Range("E39").Value = "  -0.22%  "

$ws.Range("E40").Value = "  -1.33%  "

$ws.Range("D41").Value = "'0.197"
$ws.Range("E41").Value = "  +4.16%  "

$ws.Range("D42").Value = "'59.65"
$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("E43").Value = "  -3.74%  "

$ws.Range("D44").Value = "'103.13"
$ws.Range("E44").Value = "  +6.22%  "

$ws.Range("D45").Value = "'0.473"
$ws.Range("E45").Value = "  +17.50%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.44"
$ws.Range("E46").Value = "  +11.10%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'8.27"
$ws.Range("E47").Value = "  -1.67%  "

$ws.Range("D48").Value = "'0.0971"
$ws.Range("E48").Value = "  +0.68%  "

$ws.Range("E49").Value = "  +1.04%  "

$ws.Range("D50").Value = "'1.13"
$ws.Range("E50").Value = "  +1.24%  "

$ws.Range("E51").Value = "  +0.46%  "
